# Update link for user function
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Login ---
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:1337/japtool/auth")
$ws.Range("C6").Value = "Sign in"
$ws.Range("D6").Value = "Đăng nhập"
$ws.Range("E6").Value = "Y"

# --- Row 7: Signup ---
$ws.Range("C7").Value = "Sign up"
$ws.Range("D7").Value = "Đăng ký tài khoản mới"
$ws.Range("E7").Value = "Y"

# --- Row 8: Forgot password ---
$ws.Range("B8").Value = "http://localhost:1337/japtool/forgotPassword.html"
$ws.Range("C8").Value = "Forgot password"
$ws.Range("D8").Value = "Lấy lại mật khẩu"
$ws.Range("E8").Value = "N"

# --- Row 9: User profile ---
$ws.Hyperlinks.Add($ws.Range("B9"), "http://localhost:1337/japtool/user/show/userID")
$ws.Range("C9").Value = "User Profile"
$ws.Range("D9").Value = "View user profile"
$ws.Range("E9").Value = "Y"

# --- Row 10: Edit profile ---
$ws.Range("D10").Value = "Edit user profile"
$ws.Range("E10").Value = "Y"

# --- Row 11: Change password ---
$ws.Range("D11").Value = "Change password"
$ws.Range("E11").Value = "Y"

# --- Row 12: Change avatar ---
$ws.Range("D12").Value = "Change avatar"
$ws.Range("E12").Value = "Y"

# --- Row 13: Friends function ---
$ws.Range("D13").Value = "Friends function"
$ws.Range("E13").Value = "N"

# --- Column widths (B got wider, C got narrower) ---
$ws.Columns.Item(2).ColumnWidth = 46.7
$ws.Columns.Item(3).ColumnWidth = 14.8

# --- Page setup: zoom/scale changed from 85 to 66 ---
$ws.PageSetup.Zoom = 66

# --- Selection moved to E15 ---
$ws.Range("E15").Select()
